$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals_G = @(14.04383033333333,14.04383033333333,14.04383033333333,14.04383033333333,27.64911833333333,27.64911833333333,27.64911833333333,27.64911833333333,25.89747433333334,25.89747433333334,25.89747433333334,25.89747433333334,20.774284,20.774284,20.774284,20.774284)
for ($i = 0; $i -lt $vals_G.Length; $i++) {
    $ws.Range("G" + (2 + $i)).Value = $vals_G[$i]
}

$vals_H = @(42.131491,42.131491,42.131491,42.131491,82.947355,82.947355,82.947355,82.947355,77.692423,77.692423,77.692423,77.692423,62.322852,62.322852,62.322852,62.322852)
for ($i = 0; $i -lt $vals_H.Length; $i++) {
    $ws.Range("H" + (2 + $i)).Value = $vals_H[$i]
}

$vals_I = @(0.158930310642385,0.158930310642385,0.158930310642385,0.158930310642385,0.3128977537755355,0.3128977537755355,0.3128977537755355,0.3128977537755355,0.2930748622675039,0.2930748622675039,0.2930748622675039,0.2930748622675039,0.2350970733145757,0.2350970733145757,0.2350970733145757,0.2350970733145757)
for ($i = 0; $i -lt $vals_I.Length; $i++) {
    $ws.Range("I" + (2 + $i)).Value = $vals_I[$i]
}

$vals_J = @(0.158930310642385,0.158930310642385,0.158930310642385,0.158930310642385,0.3128977537755354,0.3128977537755354,0.3128977537755354,0.3128977537755354,0.2930748622675038,0.2930748622675038,0.2930748622675038,0.2930748622675038,0.2350970733145757,0.2350970733145757,0.2350970733145757,0.2350970733145757)
for ($i = 0; $i -lt $vals_J.Length; $i++) {
    $ws.Range("J" + (2 + $i)).Value = $vals_J[$i]
}

$vals_M = @(25.69910333333333,0.1622346666666667,46.92720933333334,220.538579,25.69910333333333,0.1622346666666667,46.92720933333334,220.538579,25.69910333333333,0.1622346666666667,46.92720933333334,220.538579,25.69910333333333,0.1622346666666667,46.92720933333334,220.538579)
for ($i = 0; $i -lt $vals_M.Length; $i++) {
    $ws.Range("M" + (2 + $i)).Value = $vals_M[$i]
}

$vals_N = @(77.09731,0.486704,140.781628,661.615737,77.09731,0.486704,140.781628,661.615737,77.09731,0.486704,140.781628,661.615737,77.09731,0.486704,140.781628,661.615737)
for ($i = 0; $i -lt $vals_N.Length; $i++) {
    $ws.Range("N" + (2 + $i)).Value = $vals_N[$i]
}

$vals_O = @(0.08761243344445813,0.0005530844306649811,0.1599825079935015,0.7518519741313753,0.08761243344445813,0.0005530844306649811,0.1599825079935015,0.7518519741313753,0.08761243344445813,0.0005530844306649811,0.1599825079935015,0.7518519741313753,0.08761243344445813,0.0005530844306649811,0.1599825079935015,0.7518519741313753)
for ($i = 0; $i -lt $vals_O.Length; $i++) {
    $ws.Range("O" + (2 + $i)).Value = $vals_O[$i]
}

$vals_P = @(0.08761243344445814,0.0005530844306649812,0.1599825079935016,0.7518519741313754,0.08761243344445814,0.0005530844306649812,0.1599825079935016,0.7518519741313754,0.08761243344445814,0.0005530844306649812,0.1599825079935016,0.7518519741313754,0.08761243344445814,0.0005530844306649812,0.1599825079935016,0.7518519741313754)
for ($i = 0; $i -lt $vals_P.Length; $i++) {
    $ws.Range("P" + (2 + $i)).Value = $vals_P[$i]
}

$vals_Q = @(360.9138469321344,2.278396132851555,659.0377658941497,3097.206385430429,710.5575491238943,4.485645496435556,1297.495963910438,6097.697267836182,665.541868964681,4.201468115976889,1215.296199244961,5711.392189162306,533.8804711920133,3.370309039978666,974.8791740181173,4581.531073102436)
for ($i = 0; $i -lt $vals_Q.Length; $i++) {
    $ws.Range("Q" + (2 + $i)).Value = $vals_Q[$i]
}

$vals_R = @(3248.224622389209,20.505565195664,5931.339893047348,27874.85746887386,6395.01794211505,40.37080946792,11677.46367519394,54879.27541052563,5989.87682068213,37.813213043792,10937.66579320465,51402.52970246075,4804.924240728119,30.332781359808,8773.912566163057,41233.77965792192)
for ($i = 0; $i -lt $vals_R.Length; $i++) {
    $ws.Range("R" + (2 + $i)).Value = $vals_R[$i]
}

$vals_S = @(0.01392427126346301,0.0000879018803770521,0.02542606969275504,0.1194920678057899,0.02741373362757955,0.0001730588760032935,0.05005816739454327,0.2352527938774093,0.02567700186465542,0.0001620951433394402,0.0468868514954053,0.2203489137641037,0.02059742668876016,0.0001300285309451954,0.03761141941079792,0.1767581986840724)
for ($i = 0; $i -lt $vals_S.Length; $i++) {
    $ws.Range("S" + (2 + $i)).Value = $vals_S[$i]
}

$vals_T = @(0.01392427126346301,0.0000879018803770521,0.02542606969275504,0.1194920678057899,0.02741373362757955,0.0001730588760032935,0.05005816739454327,0.2352527938774093,0.02567700186465542,0.0001620951433394402,0.0468868514954053,0.2203489137641037,0.02059742668876016,0.0001300285309451954,0.03761141941079792,0.1767581986840724)
for ($i = 0; $i -lt $vals_T.Length; $i++) {
    $ws.Range("T" + (2 + $i)).Value = $vals_T[$i]
}


Write-Host "Edit complete"
